$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date header in C1
$ws.Range("C1").Value = "19_02_22"

# Capture current A4:B4 (A10 / Jimmy Kimmel) before shifting rows up
$rollA4 = $ws.Range("A4").Value2
$nameB4 = $ws.Range("B4").Value2

# Shift rows 5-12 up into rows 4-11 (A and B columns)
for ($r = 5; $r -le 12; $r++) {
    $ws.Cells.Item($r - 1, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r - 1, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# Place the original row 4 content (A10 / Jimmy Kimmel) at the bottom, row 12
$ws.Range("A12").Value = $rollA4
$ws.Range("B12").Value = $nameB4

# Clear row 2 completely (A2:C2)
$ws.Range("A2:C2").Clear()

# Set the new C column values of 1 for rows 7 through 12
for ($r = 7; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}
